$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.816.30'
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").Value = '1.631.21'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("D4").Value = '''0.996'
$ws.Range("E4").Value = '  -0.67%  '
$ws.Range("D5").Value = '''214.18'
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  -0.63%  '
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("E9").Value = '  -0.59%  '
$ws.Range("D10").Value = '''19.68'
$ws.Range("E10").Value = '  +0.70%  '
$ws.Range("D11").Value = '''0.0789'
$ws.Range("E11").Value = '  +0.78%  '
$ws.Range("E12").Value = '  +0.20%  '
$ws.Range("D13").Value = '1.855.89'
$ws.Range("D14").Value = '1.630.66'
$ws.Range("E14").Value = '  +0.12%  '
$ws.Range("D15").Value = '''0.553'
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("E16").Value = '  -0.18%  '
$ws.Range("D17").Value = '''62.77'
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").Value = '25.803.15'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("D19").Value = '''0.997'
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("D20").Value = '''4.44'
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D21").Value = '''191.26'
$ws.Range("E21").Value = '  -1.34%  '
$ws.Range("D22").Value = '''9.93'
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").Value = '''6.28'
$ws.Range("E23").Value = '  +0.88%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = '''1.81'
$ws.Range("E24").Value = '  +1.99%  '
$ws.Range("B25").Value = 'BinanceUSD'
$ws.Range("C25").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D25").Value = '''0.996'
$ws.Range("E25").Value = '  -0.77%  '
$ws.Range("D26").Value = '''142.34'
$ws.Range("E26").Value = '  +1.94%  '
$ws.Range("D27").Value = '''0.124'
$ws.Range("E27").Value = '  +2.78%  '
$ws.Range("D28").Value = '''6.83'
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").Value = '''15.51'
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("D30").Value = '''1.23'
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("D31").Value = '''0.0495'
$ws.Range("E31").Value = '  +1.45%  '
$ws.Range("E32").Value = '  +0.26%  '
$ws.Range("D33").Value = '''3.22'
$ws.Range("E33").Value = '  -0.56%  '
$ws.Range("D34").Value = '''1.59'
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").Value = '''0.906'
$ws.Range("E36").Value = '  +1.13%  '
$ws.Range("D37").Value = '1.142.89'
$ws.Range("E37").Value = '  +3.16%  '
$ws.Range("E38").Value = '  -0.24%  '
$ws.Range("E39").Value = '  -2.18%  '
$ws.Range("E40").Value = '  +0.17%  '
$ws.Range("D41").Value = '''0.995'
$ws.Range("E41").Value = '  -0.69%  '
$ws.Range("D43").Value = '''5.58'
$ws.Range("E43").Value = '  +0.51%  '
$ws.Range("D44").Value = '''100.65'
$ws.Range("E44").Value = '  +0.71%  '
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("D46").Value = '1.766.10'
$ws.Range("E46").Value = '  +0.41%  '
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D47").Value = '0.0₆0109'
$ws.Range("E47").Value = '  -1.70%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = '''55.36'
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '''1.48'
$ws.Range("E49").Value = '  +7.48%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.0512'
$ws.Range("E50").Value = '  +2.19%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '''0.416'
$ws.Range("E51").Value = '  -0.37%  '
